$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing "analysis method independent" / "analysis method undefined"
# values in column A down by one row to make room for the new entry, then
# insert the new "statistical model checking" value.
$ws.Range("A23").Value = $ws.Range("A22").Value2
$ws.Range("A22").Value = $ws.Range("A21").Value2
$ws.Range("A21").Value = $null
$ws.Range("A20").Value = "statistical model checking"

# Update the view state to match the recorded selection/scroll position.
$ws.Range("A22:A23").Select()
$excel.ActiveWindow.ScrollRow = 7
